$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# ---- Header row (row 1) ----
$headers = @(
    "Date",
    "Model Name",
    "Exact Precision (Micro Avg)",
    "Exact Recall (Micro Avg)",
    "Exact F1 Score (Micro Avg)",
    "Exact Precision (Macro Avg)",
    "Exact Recall (Macro Avg)",
    "Exact F1 Score (Macro Avg)",
    "Exact Precision (Weighted Avg)",
    "Exact Recall (Weighted Avg)",
    "Exact F1 Score (Weighted Avg)",
    "Partial Precision",
    "Partial Recall",
    "Partial F1 Score",
    "Partial TP",
    "Partial FP",
    "Partial FN",
    "Support",
    "Accuracy",
    "Result Link",
    "Stats Link",
    "No of GPU Used",
    "Power Consumption"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Style the header row: bold font, thin box border, centered horizontally, top-aligned vertically
# Build the style on a single cell first, then fan it out via copy/paste-special so we
# don't leave behind a trail of unused intermediate cell-format records.
$styleCell = $ws.Range("A1")
$styleCell.Font.Bold = $true
$styleCell.Borders.LineStyle = 1
$styleCell.HorizontalAlignment = -4108
$styleCell.VerticalAlignment = -4160
$styleCell.Copy()
$ws.Range("B1:W1").PasteSpecial(-4122)

# ---- Data row (row 2) ----
# Column A holds a date-like string; force it to stay text instead of being parsed as a date
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "09/10/2025"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "Qwen2.5-72B-Instruct"

$ws.Range("C2").Value = 0.3951612903225806
$ws.Range("D2").Value = 0.32996632996633
$ws.Range("E2").Value = 0.3596330275229358
$ws.Range("F2").Value = 0.1974136577708006
$ws.Range("G2").Value = 0.1450739794864496
$ws.Range("H2").Value = 0.1640170934998521
$ws.Range("I2").Value = 0.4669142145332621
$ws.Range("J2").Value = 0.32996632996633
$ws.Range("K2").Value = 0.3814513122768086
$ws.Range("L2").Value = 0.4979591836734694
$ws.Range("M2").Value = 0.4121621621621622
$ws.Range("N2").Value = 0.4510166358595194
$ws.Range("O2").Value = 122
$ws.Range("P2").Value = 123
$ws.Range("Q2").Value = 174
$ws.Range("R2").Value = 297
$ws.Range("S2").Value = 0.9451168364289994

$ws.Range("T2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-72B-Instruct_3_shot.txt"
$ws.Range("U2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-72B-Instruct_3_shot.txt"
$ws.Range("V2").Value = "4 MLGPU"
$ws.Range("W2").Value = "0.145 kWh"
$ws.Range("X2").Value = 4957
